$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.322.23"
$ws.Range("E2").Value = "  +0.64%  "

$ws.Range("D3").Value = "1.873.80"
$ws.Range("E3").Value = "  +0.62%  "

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").Value = "0.7125"
$ws.Range("E5").Value = "  +0.38%  "

$ws.Range("E6").Value = "  +0.15%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("E8").Value = "  +0.44%  "

$ws.Range("D9").Value = "0.07770"
$ws.Range("E9").Value = "  +1.91%  "

$ws.Range("D10").Value = "25.11"
$ws.Range("E10").Value = "  +1.96%  "

$ws.Range("D11").Value = "0.08405"
$ws.Range("E11").Value = "  +0.65%  "

$ws.Range("D12").Value = "1.870.28"
$ws.Range("E12").Value = "  +0.31%  "

$ws.Range("D13").Value = "5.247"
$ws.Range("E13").Value = "  +0.92%  "

$ws.Range("D14").Value = "0.7119"
$ws.Range("E14").Value = "  +0.67%  "

$ws.Range("D15").Value = "91.17"
$ws.Range("E15").Value = "  +0.10%  "

$ws.Range("D16").Value = "29.329.74"
$ws.Range("E16").Value = "  +0.48%  "

$ws.Range("D17").Value = "6.065"
$ws.Range("E17").Value = "  +2.78%  "

$ws.Range("D18").Value = "0.000008185"
$ws.Range("E18").Value = "  +4.94%  "

$ws.Range("D19").Value = "239.85"
$ws.Range("E19").Value = "  -1.13%  "

$ws.Range("E20").Value = "  +1.09%  "

$ws.Range("D21").Value = "2.118.10"
$ws.Range("E21").Value = "  +0.31%  "

$ws.Range("D22").Value = "0.9996"
$ws.Range("E22").Value = "  +0.00%  "

$ws.Range("D23").Value = "7.765"
$ws.Range("E23").Value = "  -1.36%  "

$ws.Range("D24").Value = "1.003"
$ws.Range("E24").Value = "  +0.35%  "

$ws.Range("D25").Value = "0.1587"
$ws.Range("E25").Value = "  +0.10%  "

$ws.Range("D26").Value = "162.99"
$ws.Range("E26").Value = "  -0.59%  "

$ws.Range("D27").Value = "9.031"
$ws.Range("E27").Value = "  +0.96%  "

$ws.Range("E28").Value = "  +0.48%  "

$ws.Range("E29").Value = "  +0.68%  "

$ws.Range("D30").Value = "4.405"
$ws.Range("E30").Value = "  +0.40%  "

$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").Value = "1.290"
$ws.Range("E31").Value = "  -2.27%  "

$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "4.323"
$ws.Range("E32").Value = "  +1.80%  "

$ws.Range("D33").Value = "0.05297"
$ws.Range("E33").Value = "  +3.11%  "

$ws.Range("D34").Value = "1.939"
$ws.Range("E34").Value = "  +1.54%  "

$ws.Range("E35").Value = "  +1.27%  "

$ws.Range("D36").Value = "0.7436"
$ws.Range("E36").Value = "  -6.36%  "

$ws.Range("D37").Value = "2.705"
$ws.Range("E37").Value = "  +0.69%  "

$ws.Range("D38").Value = "0.01876"
$ws.Range("E38").Value = "  +1.73%  "

$ws.Range("D39").Value = "1.222.99"
$ws.Range("E39").Value = "  +5.15%  "

$ws.Range("D40").Value = "2.727"
$ws.Range("E40").Value = "  +1.15%  "

$ws.Range("D41").Value = "6.542"
$ws.Range("E41").Value = "  +5.05%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "0.8858"
$ws.Range("E42").Value = "  -0.59%  "

$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "109.44"
$ws.Range("E43").Value = "  +6.61%  "

$ws.Range("D44").Value = "72.49"

$ws.Range("E45").Value = "  +0.18%  "

$ws.Range("D46").Value = "2.017.13"
$ws.Range("E46").Value = "  +0.35%  "

$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Value = "1.796"
$ws.Range("E47").Value = "  +1.21%  "

$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "0.5196"
$ws.Range("E48").Value = "  +0.20%  "

$ws.Range("D49").Value = "0.00000000123"
$ws.Range("E49").Value = "  +5.16%  "

$ws.Range("D50").Value = "9.374"
$ws.Range("E50").Value = "  +1.09%  "

$ws.Range("E51").Value = "  +0.97%  "
